$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 7
$ws.Range("B3").Value = 3
$ws.Range("B5").Value = 3
$ws.Range("C8").Value = 71
$ws.Range("C9").Value = 3
$ws.Range("C10").Value = 2
$ws.Range("C11").Value = 2
$ws.Range("C12").Value = 2
$ws.Range("C13").Value = 1
$ws.Range("B14").Value = 6
$ws.Range("B15").Value = 4
$ws.Range("B16").Value = 4
$ws.Range("B17").Value = 3
$ws.Range("B18").Value = 4
$ws.Range("B19").Value = 3
$ws.Range("C20").Value = 49
$ws.Range("C21").Value = 4
$ws.Range("C22").Value = 1085
$ws.Range("C23").Value = 3
$ws.Range("C24").Value = 2
$ws.Range("C25").Value = 940
